# Updated cryptos list on Thu Feb  8 18:00:40 UTC 2024 with GitHub Actions
# Refreshes Price (column D) and Volume(1h) (column E) for each coin row,
# and reflects that Stacks moved ahead of FraxShare in the ranking
# (their B/C/D/E cell contents swap rows 48 and 49).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.062.48'
$ws.Range("E2").Value = '  +3.54%  '
$ws.Range("D3").Value = '2.426.05'
$ws.Range("E3").Value = '  +0.60%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''317.13'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.61%  '
$ws.Range("D6").Value = '''102.59'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.33%  '
$ws.Range("E7").Value = '  +1.17%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").Value = '''0.525'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +7.27%  '
$ws.Range("D10").Value = '''35.39'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.03%  '
$ws.Range("E11").Value = '  +0.67%  '
$ws.Range("E12").Value = '  -2.53%  '
$ws.Range("D13").Value = '''18.13'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.86%  '
$ws.Range("D14").Value = '''6.99'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.53%  '
$ws.Range("D15").Value = '2.808.46'
$ws.Range("D16").Value = '2.425.84'
$ws.Range("E16").Value = '  -1.99%  '
$ws.Range("D17").Value = '''0.836'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.34%  '
$ws.Range("D18").Value = '45.044.23'
$ws.Range("E18").Value = '  +3.58%  '
$ws.Range("D19").Value = '''12.24'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.66%  '
$ws.Range("D20").Value = '''6.36'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.87%  '
$ws.Range("E21").Value = '  +2.19%  '
$ws.Range("E22").Value = '  +0.51%  '
$ws.Range("D23").Value = '''243.89'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.69%  '
$ws.Range("E24").Value = '  +1.09%  '
$ws.Range("D25").Value = '''2.50'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.16%  '
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("D27").Value = '''25.27'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.44%  '
$ws.Range("D28").Value = '''2.18'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.45%  '
$ws.Range("E29").Value = '  +1.38%  '
$ws.Range("D30").Value = '''49.21'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.02%  '
$ws.Range("D31").Value = '''32.73'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.14%  '
$ws.Range("D32").Value = '''20.20'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +10.32%  '
$ws.Range("E33").Value = '  +9.61%  '
$ws.Range("E34").Value = '  +1.61%  '
$ws.Range("E35").Value = '  +0.31%  '
$ws.Range("E36").Value = '  +2.46%  '
$ws.Range("E37").Value = '  -0.80%  '
$ws.Range("D38").Value = '''4.41'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.45%  '
$ws.Range("E39").Value = '  -2.54%  '
$ws.Range("D40").Value = '''125.23'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -7.05%  '
$ws.Range("E41").Value = '  -2.31%  '
$ws.Range("E42").Value = '  +0.57%  '
$ws.Range("D43").Value = '''20.80'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.32%  '
$ws.Range("D44").Value = '''0.0288'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.66%  '
$ws.Range("D45").Value = '1.932.93'
$ws.Range("E45").Value = '  -0.93%  '
$ws.Range("E46").Value = '  -2.48%  '
$ws.Range("E47").Value = '  +3.80%  '
$ws.Range("D50").Value = '''76.35'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.74%  '
$ws.Range("E51").Value = '  +2.58%  '

# Row 48/49: FraxShare and Stacks swap positions with updated values
$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").Value = '''1.82'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +17.10%  '
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").Value = '''9.24'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.45%  '
